$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.921780705451965
$ws.Range("B1").Value = 2.798107862472534
$ws.Range("C1").Value = 3.078586578369141
$ws.Range("D1").Value = 2.636105298995972
$ws.Range("E1").Value = 1.022797584533691
